$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B2 text
$ws.Range("B2").Value = "LIVEHTA_723 - Test723 - 1/13/2023"

# Rebuild column I (I2:I10) explicitly, since a "Report-" row is inserted
# at I2 pushing the existing entries down, and new Economic/Quality of
# Life/Real-world Evidence report rows are appended at the end.
$ws.Range("I2").Value = "Report-"
$ws.Range("I3").Value = "ExcelReport-LIVEHTA_723 - Test723-Clinical-"
$ws.Range("I4").Value = "WordReport-LIVEHTA_723 - Test723-Clinical-"
$ws.Range("I5").Value = "ExcelReport-LIVEHTA_723 - Test723-Economic-"
$ws.Range("I6").Value = "WordReport-LIVEHTA_723 - Test723-Economic-"
$ws.Range("I7").Value = "ExcelReport-LIVEHTA_723 - Test723-Quality of Life-"
$ws.Range("I8").Value = "WordReport-LIVEHTA_723 - Test723-Quality of Life-"
$ws.Range("I9").Value = "ExcelReport-LIVEHTA_723 - Test723-Real-world Evidence-"
$ws.Range("I10").Value = "WordReport-LIVEHTA_723 - Test723-Real-world Evidence-"

$ws.Columns.Item(9).ColumnWidth = 48.67

# Force a pageSetup element with portrait orientation (matches the page
# setup block that Excel wrote when the print settings were touched).
$ws.PageSetup.Orientation = 1

$ws.Range("B2").Select()
